# The "Estado de Cuenta" worker/period table (rows 16-26) is being
# re-sorted: previously grouped by worker (each worker's periods listed
# 2406 -> 2403 descending), the table is now grouped by period ascending
# (2403/2404 first, then 2405 for everyone, then 2406 for everyone), with
# the underlying (worker, periodo, valor) data unchanged.
#
# Source row -> destination row mapping:
#   19 -> 16   18 -> 17   22 -> 18   25 -> 19   17 -> 20
#   21 -> 21   24 -> 22   16 -> 23   20 -> 24   23 -> 25   26 -> 26

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 16
$lastRow = 26

# Snapshot all current values for columns B..G across the affected rows
# before writing anything back, so reads are never polluted by writes.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapshot[$r] = @{
        B = $ws.Cells.Item($r, 2).Value()
        C = $ws.Cells.Item($r, 3).Value()
        D = $ws.Cells.Item($r, 4).Value()
        E = $ws.Cells.Item($r, 5).Value()
        F = $ws.Cells.Item($r, 6).Value()
        G = $ws.Cells.Item($r, 7).Value()
    }
}

$rowMap = @{
    16 = 19
    17 = 18
    18 = 22
    19 = 25
    20 = 17
    21 = 21
    22 = 24
    23 = 16
    24 = 20
    25 = 23
    26 = 26
}

foreach ($destRow in ($rowMap.Keys | Sort-Object)) {
    $srcRow = $rowMap[$destRow]
    $vals = $snapshot[$srcRow]

    $ws.Cells.Item($destRow, 2).Value = $vals.B
    $ws.Cells.Item($destRow, 3).Value = $vals.C
    $ws.Cells.Item($destRow, 4).Value = $vals.D
    $ws.Cells.Item($destRow, 5).Value = $vals.E
    $ws.Cells.Item($destRow, 6).Value = $vals.F
    $ws.Cells.Item($destRow, 7).Value = $vals.G
}
